$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = "Georgetown"
$ws.Range("B17").Value = "Hoyas"
$ws.Range("C17").Value = "MAAC"
$ws.Range("D17").Value = "NotAName#0591"
$ws.Range("E17").Value = "Nota Nam"
$ws.Range("F17").Value = "west coast"
$ws.Range("G17").Value = "4-3"

$ws.Range("A18").Value = "South East North Western Wyoming A&M Tech State"
$ws.Range("B18").Value = "Country Girls That Make Do"
$ws.Range("C18").Value = "ECFL"
$ws.Range("D18").Value = "PleaseEndMeNow#3186"
$ws.Range("E18").Value = "Brett Hundley"
$ws.Range("F18").Value = "spread"
$ws.Range("G18").Value = "4-3"

$ws.Range("A19").Value = "Nebraska"
$ws.Range("B19").Value = "Nebraskans"
$ws.Range("C19").Value = "The Wonderful Nebraska Coalition"
$ws.Range("D19").Value = "Hobbes .T. Hero#4989"
$ws.Range("E19").Value = "Paris Riley"
$ws.Range("F19").Value = "air raid"
$ws.Range("G19").Value = "4-3"

$ws.Range("A20").Value = "Liberty"
$ws.Range("B20").Value = "Cool Guys"
$ws.Range("C20").Value = "ACC"
$ws.Range("D20").Value = "jakeysnakey#6969"
$ws.Range("E20").Value = "Jakob"
$ws.Range("F20").Value = "air raid"
$ws.Range("G20").Value = "5-2"

$ws.Range("A21").Value = "Milk University"
$ws.Range("B21").Value = "Milkmen"
$ws.Range("C21").Value = "Milk 12"
$ws.Range("D21").Value = "Naki#2555"
$ws.Range("E21").Value = 'Tony "The Milkman" Stevens'
$ws.Range("F21").Value = "flexbone"
$ws.Range("G21").Value = "4-4"

$ws.Range("A22").Value = "Florida"
$ws.Range("B22").Value = "Footballers"
$ws.Range("C22").Value = "Floridan Football Fantasy"
$ws.Range("D22").Value = "JVitt#8369"
$ws.Range("E22").Value = "JVitt"
$ws.Range("F22").Value = "air raid"
$ws.Range("G22").Value = "4-4"

$ws.Range("A23").Value = "Cudahy"
$ws.Range("B23").Value = "Chuds"
$ws.Range("C23").Value = "The Resistance"
$ws.Range("D23").Value = "penguino#2114"
$ws.Range("E23").Value = "Pingu"
$ws.Range("F23").Value = "air raid"
$ws.Range("G23").Value = "5-2"

$ws.Range("A24").Value = "North Atlanta"
$ws.Range("B24").Value = "Hornets"
$ws.Range("C24").Value = "C-USA"
$ws.Range("D24").Value = "Starboy#1512"
$ws.Range("E24").Value = "Jeff Hollins"
$ws.Range("F24").Value = "west coast"
$ws.Range("G24").Value = "3-4"

$ws.Range("A25").Value = "Mommy’s"
$ws.Range("B25").Value = "Milkies"
$ws.Range("C25").Value = "JUG"
$ws.Range("D25").Value = "lancer52#4833"
$ws.Range("E25").Value = "Dick Sux"
$ws.Range("F25").Value = "flexbone"
$ws.Range("G25").Value = "5-2"

$ws.Range("A26").Value = "USC"
$ws.Range("B26").Value = "Trojans"
$ws.Range("C26").Value = "PAC 12"
$ws.Range("D26").Value = "stinkywrestler#7847"
$ws.Range("E26").Value = "Oliver Raymond"
$ws.Range("F26").Value = "spread"
$ws.Range("G26").Value = "4-3"

$ws.Range("A27").Value = "Gushbaba"
$ws.Range("B27").Value = "Gushbabenbabens"
$ws.Range("C27").Value = "The Wonderful Nebraska Coalition"
$ws.Range("D27").Value = "Pizza Chef#2639"
$ws.Range("E27").Value = "Gushbab"
$ws.Range("F27").Value = "air raid"
